# Power Storage sheet: set ExisUnits (col E) to 0 for rows 7 and 10,
# and MaxInvest / "MaxlineLoad" (col S) to 200 for rows 7-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

$ws.Range("E7").Value = 0
$ws.Range("E10").Value = 0

$ws.Range("S7:S11").Value = 200
